$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 53
$prev = $row - 1

# Copy the formatting of the row above (same table), then overwrite the
# values cell-by-cell so the new row matches the rest of the sheet's style
# (bold/bordered/centered "Indice" column, and the datetime column).
$ws.Range("A" + $prev + ":V" + $prev).Copy()
$ws.Range("A" + $row + ":V" + $row).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item($row, 1).Value = 52
$ws.Cells.Item($row, 2).Value = "wales"
$ws.Cells.Item($row, 3).Value = "cymru-premier"
$ws.Cells.Item($row, 4).Value = "2023-2024"
$ws.Cells.Item($row, 5).Value = 45196.86458333334
$ws.Cells.Item($row, 6).Value = "Bala"
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = "Aberystwyth"
$ws.Cells.Item($row, 9).Value = 1
$ws.Cells.Item($row, 10).Value = 1.37
$ws.Cells.Item($row, 11).Value = "26/09/2023 08:13"
$ws.Cells.Item($row, 12).Value = 1.3
$ws.Cells.Item($row, 13).Value = "27/09/2023 14:42"
$ws.Cells.Item($row, 14).Value = 4.7
$ws.Cells.Item($row, 15).Value = "26/09/2023 08:13"
$ws.Cells.Item($row, 16).Value = 5.61
$ws.Cells.Item($row, 17).Value = "27/09/2023 20:29"
$ws.Cells.Item($row, 18).Value = 6.61
$ws.Cells.Item($row, 19).Value = "26/09/2023 08:13"
$ws.Cells.Item($row, 20).Value = 9.23
$ws.Cells.Item($row, 21).Value = "27/09/2023 20:29"
$ws.Cells.Item($row, 22).Value = "https://www.betexplorer.com/football/wales/cymru-premier/bala-aberystwyth/IF1445Zt/"
